$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws 'D2' '256.56'
Set-TextValue $ws 'E2' '0.53%'
Set-TextValue $ws 'D3' '27.01'
Set-TextValue $ws 'E3' '-4.08%'
Set-TextValue $ws 'D4' '4.719'
Set-TextValue $ws 'E4' '-10.17%'
Set-TextValue $ws 'D5' '0.05944'
Set-TextValue $ws 'E5' '1.58%'
Set-TextValue $ws 'D6' '6.659'
Set-TextValue $ws 'E6' '-0.69%'
Set-TextValue $ws 'D7' '0.8699'
Set-TextValue $ws 'E7' '0.21%'
Set-TextValue $ws 'D8' '0.9510'
Set-TextValue $ws 'E8' '-8.26%'
Set-TextValue $ws 'E9' '-0.73%'
Set-TextValue $ws 'D10' '0.03932'
Set-TextValue $ws 'E10' '13.15%'
Set-TextValue $ws 'D11' '0.07160'
Set-TextValue $ws 'E11' '0.14%'
Set-TextValue $ws 'D12' '0.03193'
Set-TextValue $ws 'E12' '0.21%'
Set-TextValue $ws 'D13' '0.09242'
Set-TextValue $ws 'E13' '0.18%'
Set-TextValue $ws 'D14' '0.001547'
Set-TextValue $ws 'E14' '-0.11%'
Set-TextValue $ws 'D15' '0.0006049'
Set-TextValue $ws 'E15' '-0.54%'
Set-TextValue $ws 'D16' '0.006052'
Set-TextValue $ws 'E16' '3.86%'
Set-TextValue $ws 'D17' '3.484'
Set-TextValue $ws 'E17' '-0.41%'
Set-TextValue $ws 'D18' '3.202'
Set-TextValue $ws 'E18' '-0.95%'
Set-TextValue $ws 'E19' '0.76%'
Set-TextValue $ws 'D20' '0.3134'
Set-TextValue $ws 'E20' '-1.45%'
Set-TextValue $ws 'E21' '-1.41%'
Set-TextValue $ws 'D22' '3.810'
Set-TextValue $ws 'E22' '7.01%'
Set-TextValue $ws 'D23' '0.04203'
Set-TextValue $ws 'E23' '0.86%'
Set-TextValue $ws 'E24' '2.44%'
Set-TextValue $ws 'D25' '0.001219'
Set-TextValue $ws 'E25' '-0.53%'
Set-TextValue $ws 'D26' '0.004497'
Set-TextValue $ws 'E26' '-7.84%'
Set-TextValue $ws 'D28' '0.0001936'
Set-TextValue $ws 'E28' '142.01%'
Set-TextValue $ws 'D40' '0.03830'
Set-TextValue $ws 'E40' '0.23%'
$ws.Range('B41').Value = 'BKEXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue $ws 'D41' '0.1100'
Set-TextValue $ws 'E41' '-0.32%'
$ws.Range('B42').Value = 'KickToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue $ws 'D42' '0.003968'
Set-TextValue $ws 'E42' '-31.73%'
Set-TextValue $ws 'E43' '-3.74%'
Set-TextValue $ws 'D44' '0.01058'
Set-TextValue $ws 'E44' '9.27%'
Set-TextValue $ws 'D45' '0.00005498'
Set-TextValue $ws 'E45' '5.11%'
Set-TextValue $ws 'E46' '-0.05%'
Set-TextValue $ws 'D47' '0.08850'
Set-TextValue $ws 'E47' '-4.86%'
Set-TextValue $ws 'D48' '0.002393'
Set-TextValue $ws 'E48' '11.20%'
Set-TextValue $ws 'D49' '0.00002099'
Set-TextValue $ws 'E49' '-0.05%'
Set-TextValue $ws 'D50' '0.0001999'
Set-TextValue $ws 'E50' '-0.05%'
